$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 15, columns B:AW from 30 to 0.5
$ws.Range("B15:AW15").Value = 0.5

# Update the selection on the active sheet
$ws.Activate()
$ws.Range("AS21").Select()
